$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.483.65'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.924.11'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '375.26'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +6.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.52'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.541'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.58%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.584'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.78'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.74%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0838'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.30'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.386.28'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.35'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.923.65'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.937'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -6.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.422.52'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.39'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.95'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0946'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.31'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '261.45'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.77'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.66%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.12'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.97%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.167'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.23%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.34'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.78'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.99'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +9.96%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -6.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.81'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.82'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.36%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '33.93'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.31%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.01'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -8.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.91'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.60'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -7.36%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.49%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '124.63'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.85'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.06'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.270'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +12.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.018.10'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.87%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.95%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.210.85'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.64%  '
